$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Add a new date-stamped cell in C1, using the built-in "mm-dd-yy" (numFmtId 14)
# date number format, then assign the date value (2021-04-21 -> serial 44307).
$ws.Range("C1").NumberFormat = "mm-dd-yy"
$ws.Range("C1").Value = (Get-Date -Year 2021 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0)
